# Update timestamps related to the "6a72a293-1cef-46b6-89c7-c0b0778ecde2.md" row
# to reflect a freshly (re)generated handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 6a72a293 row (row 4)
$overview.Range("G4").Value = "2016-09-04 04:50:10"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for row 4
$zhcn.Range("H4").Value = "2016-09-04 04:49:58"
$zhcn.Range("K4").Value = "2016-09-04 04:50:30"

# de-de sheet: "Correspond Handoff Datetime" (mirrors the Overview HO Xliff date)
# and "Correspond Handback DateTime" for row 4
$dede.Range("H4").Value = "2016-09-04 04:50:10"
$dede.Range("K4").Value = "2016-09-04 04:50:38"
